$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# The "Periodo Mora" column (E16:E28) previously listed the periods in
# descending order (2102,2101,2012,...,2002). Re-key it in ascending order
# (2002,2003,...,2012,2101,2102) as part of the database refresh.
$periodos = @("2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}

# The "Salario Basico" value that belonged to the first row now belongs to
# the last row of the refreshed block (and vice versa).
$ws.Cells.Item(16, 6).Value = 35112
$ws.Cells.Item(28, 6).Value = 25749
